$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1000
$ws.Range("I40").Value = 970.85187
$ws.Range("J40").Value = 1087.4445
$ws.Range("K40").Value = 970.85187
$ws.Range("L40").Value = 1087.4445
$ws.Range("M40").Value = -795.85187
$ws.Range("N40").Value = -1437.4445

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 33679.8
$ws.Range("J75").Value = 33679.8
$ws.Range("L75").Value = 33679.8
$ws.Range("N75").Value = -35551.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H78").Value = 33679.8
$ws.Range("J78").Value = 33679.8
$ws.Range("L78").Value = 101039.4
$ws.Range("N78").Value = -110399.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1800
$ws.Range("J121").Value = 1800
$ws.Range("L121").Value = 5400
$ws.Range("N121").Value = -8894

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 997825.6
$ws.Range("I137").Value = 1659842.6
$ws.Range("J137").Value = 4800
$ws.Range("K137").Value = 4979527.800000001
$ws.Range("L137").Value = 14400
$ws.Range("M137").Value = -4976977.800000001
$ws.Range("N137").Value = -19500

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3974.4941
$ws.Range("I32").Value = 4068.2534
$ws.Range("J32").Value = 3499
$ws.Range("K32").Value = 4068.2534
$ws.Range("L32").Value = 3499
$ws.Range("M32").Value = -3781.2534
$ws.Range("N32").Value = -4073

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2693.9092
$ws.Range("I61").Value = 1508.625
$ws.Range("J61").Value = 5854.6665
$ws.Range("K61").Value = 1508.625
$ws.Range("L61").Value = 5854.6665
$ws.Range("M61").Value = -1296.625
$ws.Range("N61").Value = -6278.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3578.7234
$ws.Range("I74").Value = 895
$ws.Range("J74").Value = 5244.483
$ws.Range("K74").Value = 895
$ws.Range("L74").Value = 5244.483
$ws.Range("M74").Value = -21
$ws.Range("N74").Value = -6992.483

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3578.7234
$ws.Range("I77").Value = 895
$ws.Range("J77").Value = 5244.483
$ws.Range("K77").Value = 4475
$ws.Range("L77").Value = 26222.415
$ws.Range("M77").Value = -107
$ws.Range("N77").Value = -34958.415

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1382355.9
$ws.Range("I132").Value = 2031358.9
$ws.Range("J132").Value = 3224.625
$ws.Range("K132").Value = 6094076.699999999
$ws.Range("L132").Value = 9673.875
$ws.Range("M132").Value = -6091546.699999999
$ws.Range("N132").Value = -14733.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2693.9092
$ws.Range("I136").Value = 1508.625
$ws.Range("J136").Value = 5854.6665
$ws.Range("K136").Value = 4525.875
$ws.Range("L136").Value = 17563.9995
$ws.Range("M136").Value = -1975.875
$ws.Range("N136").Value = -22663.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 40215.355
$ws.Range("I134").Value = 49465.92
$ws.Range("J134").Value = 1671.3334
$ws.Range("K134").Value = 148397.76
$ws.Range("L134").Value = 5014.0002
$ws.Range("M134").Value = -145862.76
$ws.Range("N134").Value = -10084.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 53265.832
$ws.Range("J140").Value = 53265.832
$ws.Range("L140").Value = 53265.832
$ws.Range("N140").Value = -63625.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1531.5454
$ws.Range("I31").Value = 918.375
$ws.Range("J31").Value = 3166.6667
$ws.Range("K31").Value = 918.375
$ws.Range("L31").Value = 3166.6667
$ws.Range("M31").Value = -623.375
$ws.Range("N31").Value = -3756.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1531.5454
$ws.Range("I34").Value = 918.375
$ws.Range("J34").Value = 3166.6667
$ws.Range("K34").Value = 918.375
$ws.Range("L34").Value = 3166.6667
$ws.Range("M34").Value = -716.375
$ws.Range("N34").Value = -3570.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1921.4894
$ws.Range("I134").Value = 2110.1316
$ws.Range("J134").Value = 1125
$ws.Range("K134").Value = 6330.3948
$ws.Range("L134").Value = 3375
$ws.Range("M134").Value = -3795.3948
$ws.Range("N134").Value = -8445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 39625
$ws.Range("J135").Value = 39625
$ws.Range("L135").Value = 39625
$ws.Range("N135").Value = -49765

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13334145
$ws.Range("I4").Value = 28572144
$ws.Range("J4").Value = 897.5
$ws.Range("K4").Value = 85716432
$ws.Range("L4").Value = 2692.5
$ws.Range("M4").Value = -85716320
$ws.Range("N4").Value = -2916.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 12756325
$ws.Range("I122").Value = 20833842
$ws.Range("K122").Value = 187504578
$ws.Range("M122").Value = -187502128

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 33367.11
$ws.Range("I137").Value = 2735.7144
$ws.Range("J137").Value = 40760.9
$ws.Range("K137").Value = 8207.143199999999
$ws.Range("L137").Value = 122282.7
$ws.Range("M137").Value = -3107.143199999999
$ws.Range("N137").Value = -132482.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2130.726
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 2157.5417
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 2157.5417
$ws.Range("M5").Value = -88
$ws.Range("N5").Value = -2381.5417

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 530.7308
$ws.Range("I107").Value = 370
$ws.Range("J107").Value = 787.9
$ws.Range("K107").Value = 370
$ws.Range("L107").Value = 787.9
$ws.Range("M107").Value = 1550
$ws.Range("N107").Value = -4627.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7046.8667
$ws.Range("I122").Value = 7715.615
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 23146.845
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -20696.845
$ws.Range("N122").Value = -13000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3055.4
$ws.Range("I132").Value = 2220.2856
$ws.Range("J132").Value = 4118.273
$ws.Range("K132").Value = 6660.8568
$ws.Range("L132").Value = 12354.819
$ws.Range("M132").Value = -4130.8568
$ws.Range("N132").Value = -17414.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5170999.5
$ws.Range("I2").Value = 1000000
$ws.Range("J2").Value = 6005199.5
$ws.Range("K2").Value = 1000000
$ws.Range("L2").Value = 6005199.5
$ws.Range("M2").Value = -999888
$ws.Range("N2").Value = -6005423.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1758.8649
$ws.Range("I7").Value = 1599.3928
$ws.Range("J7").Value = 2255
$ws.Range("K7").Value = 1599.3928
$ws.Range("L7").Value = 2255
$ws.Range("M7").Value = -1487.3928
$ws.Range("N7").Value = -2479

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1758.8649
$ws.Range("I126").Value = 1599.3928
$ws.Range("J126").Value = 2255
$ws.Range("K126").Value = 4798.178400000001
$ws.Range("L126").Value = 6765
$ws.Range("M126").Value = -2328.178400000001
$ws.Range("N126").Value = -11705

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 75249.5
$ws.Range("I132").Value = 147157.14
$ws.Range("J132").Value = 3341.8572
$ws.Range("K132").Value = 441471.42
$ws.Range("L132").Value = 10025.5716
$ws.Range("M132").Value = -438941.42
$ws.Range("N132").Value = -15085.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 17170286
$ws.Range("J2").Value = 55000000
$ws.Range("L2").Value = 55000000
$ws.Range("N2").Value = -55000224

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 11026300
$ws.Range("I4").Value = 2536000
$ws.Range("J4").Value = 16686500
$ws.Range("K4").Value = 2536000
$ws.Range("L4").Value = 16686500
$ws.Range("M4").Value = -2535887
$ws.Range("N4").Value = -16686726

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5642.44
$ws.Range("I132").Value = 7635.8667
$ws.Range("K132").Value = 22907.6001
$ws.Range("M132").Value = -20377.6001
